$d = $word.ActiveDocument

# The document has two paragraphs right after the ">>> your stuff ..." line:
#   5: "A cool quote by Dijkstra:"
#   6: ""Computer science is no more about computers than astronomy is about telescopes.""
# They get collapsed into a single paragraph that reads:
#   "Some people think earth is flat. "
# where "Some people think earth is flat." is highlighted yellow and is
# followed by a separate, unhighlighted run containing just a space.

# 1) Remove the second paragraph (the quote) entirely, including its own
#    paragraph mark, so the two paragraphs become one.
$quotePara = $d.Paragraphs(6)
$quoteRange = $quotePara.Range.Duplicate
$quoteRange.Delete()

# 2) Replace the text of the remaining paragraph (keep its paragraph mark)
#    with the new sentence plus a trailing space.
$introPara = $d.Paragraphs(5)
$newText = "Some people think earth is flat."
$textRange = $introPara.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Text = "$newText "

# 3) Highlight only the sentence (not the trailing space) in yellow. Using
#    Find/Replacement with a formatting-only replacement correctly splits
#    the run instead of touching the whole paragraph.
$para = $d.Paragraphs(5)
$searchRange = $para.Range.Duplicate
$searchRange.Find.ClearFormatting()
$searchRange.Find.Replacement.ClearFormatting()
$searchRange.Find.Replacement.Highlight = $true
$found = $searchRange.Find.Execute($newText, $false, $false, $false, $false, `
                                    $false, $true, 1, $true, $newText, 2)
